# "fix language to lowercase" -- rename the sheet tabs to the new
# "Russian - <Category>" display names, move the active/selected tab from
# "adverbs" to "conjunctions", update a couple of cell selections, and make
# sure the "verbs" sheet carries a (portrait, letter/A4-class paper) page
# setup like its siblings already do.

$wb = $excel.ActiveWorkbook

# --- Rename the worksheet tabs (lowercase -> "Russian - X") -------------
$wb.Worksheets.Item("adjectives").Name   = "Russian - Adjectives"
$wb.Worksheets.Item("adverbs").Name      = "Russian - Adverbs"
$wb.Worksheets.Item("conjunctions").Name = "Russian - Conjunctions"
$wb.Worksheets.Item("expressions").Name  = "Russian - Expressions"
$wb.Worksheets.Item("nouns").Name        = "Russian - Nouns"
$wb.Worksheets.Item("verbs").Name        = "Russian - Verbs"

# --- "Russian - Verbs": new selection + explicit page setup -------------
$wsVerbs = $wb.Worksheets.Item("Russian - Verbs")
$wsVerbs.Range("G30").Select()
$wsVerbs.PageSetup.PaperSize = 9
$wsVerbs.PageSetup.Orientation = 1

# --- "Russian - Conjunctions" becomes the active/selected tab -----------
$wsConj = $wb.Worksheets.Item("Russian - Conjunctions")
$wsConj.Activate()
$wsConj.Range("G12").Select()
